# Apply the crypto price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# forced to text format first, otherwise Excel's automatic type detection
# would silently convert them into numeric cells (losing the original
# text/string cell type used throughout this sheet).
$forceTextCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D16", "D19", "D21", "D22", "D23", "D26", "D27", "D28", "D29", "D30", "D31", "D34", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new values cell by cell.
$ws.Range('D2').Value = '36.718.37'
$ws.Range('E2').Value = '  +3.94%  '
$ws.Range('D3').Value = '1.920.00'
$ws.Range('E3').Value = '  +2.12%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').Value = '0.703'
$ws.Range('E5').Value = '  +3.23%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '249.98'
$ws.Range('E6').Value = '  +1.06%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '44.27'
$ws.Range('E8').Value = '  +1.07%  '
$ws.Range('D9').Value = '58.67'
$ws.Range('E9').Value = '  +9.52%  '
$ws.Range('D10').Value = '0.369'
$ws.Range('E10').Value = '  +3.18%  '
$ws.Range('D11').Value = '0.0767'
$ws.Range('E11').Value = '  +3.60%  '
$ws.Range('E12').Value = '  +2.56%  '
$ws.Range('D13').Value = '14.56'
$ws.Range('E13').Value = '  +7.63%  '
$ws.Range('D14').Value = '0.831'
$ws.Range('E14').Value = '  +7.89%  '
$ws.Range('D15').Value = '2.202.86'
$ws.Range('E15').Value = '  +2.35%  '
$ws.Range('D16').Value = '5.14'
$ws.Range('E16').Value = '  +4.20%  '
$ws.Range('D17').Value = '1.916.33'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('D18').Value = '36.690.99'
$ws.Range('E18').Value = '  +3.74%  '
$ws.Range('D19').Value = '74.63'
$ws.Range('E19').Value = '  +2.72%  '
$ws.Range('D20').Value = '0.0₃0864'
$ws.Range('E20').Value = '  +4.94%  '
$ws.Range('D21').Value = '251.33'
$ws.Range('E21').Value = '  +2.96%  '
$ws.Range('D22').Value = '13.40'
$ws.Range('E22').Value = '  +4.02%  '
$ws.Range('D23').Value = '5.24'
$ws.Range('E23').Value = '  +4.79%  '
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = '2.21'
$ws.Range('E26').Value = '  -2.11%  '
$ws.Range('D27').Value = '168.15'
$ws.Range('E27').Value = '  +1.30%  '
$ws.Range('D28').Value = '8.87'
$ws.Range('E28').Value = '  +3.09%  '
$ws.Range('D29').Value = '18.76'
$ws.Range('E29').Value = '  +2.51%  '
$ws.Range('D30').Value = '0.129'
$ws.Range('E30').Value = '  +1.79%  '
$ws.Range('D31').Value = '4.66'
$ws.Range('E31').Value = '  +8.46%  '
$ws.Range('E32').Value = '  +4.87%  '
$ws.Range('E33').Value = '  -3.43%  '
$ws.Range('D34').Value = '4.36'
$ws.Range('E34').Value = '  +4.44%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('E36').Value = '  -3.88%  '
$ws.Range('D37').Value = '0.0870'
$ws.Range('E37').Value = '  +18.69%  '
$ws.Range('D38').Value = '0.902'
$ws.Range('E38').Value = '  +6.27%  '
$ws.Range('D39').Value = '17.46'
$ws.Range('E39').Value = '  +47.48%  '
$ws.Range('E40').Value = '  +6.10%  '
$ws.Range('D41').Value = '106.93'
$ws.Range('E41').Value = '  +10.32%  '
$ws.Range('D42').Value = '0.0229'
$ws.Range('E42').Value = '  +4.64%  '
$ws.Range('D43').Value = '17.56'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').Value = '1.11'
$ws.Range('E44').Value = '  +3.42%  '
$ws.Range('D45').Value = '2.84'
$ws.Range('E45').Value = '  +18.89%  '
$ws.Range('D46').Value = '1.347.27'
$ws.Range('E46').Value = '  +2.62%  '
$ws.Range('D47').Value = '2.41'
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('D48').Value = '0.0815'
$ws.Range('E48').Value = '  +1.48%  '
$ws.Range('D49').Value = '2.80'
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('D50').Value = '6.43'
$ws.Range('E50').Value = '  +2.46%  '
$ws.Range('D51').Value = '43.86'
$ws.Range('E51').Value = '  +4.19%  '
